$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header: row 44 changes from "Level 5-1" to "Level 5-2" ---
$ws.Cells.Item(44, 1).Value2 = "Level 5-2"

# --- Row 58: B58 20543 -> 20537, add E58 = 20600 ---
$ws.Cells.Item(58, 2).Value2 = 20537
$ws.Cells.Item(58, 5).Value2 = 20600

# --- Row 59: B59 21058 -> 21051, add E59 = 21114 ---
$ws.Cells.Item(59, 2).Value2 = 21051
$ws.Cells.Item(59, 5).Value2 = 21114

# --- New row 60 data ---
$ws.Cells.Item(60, 1).Value2 = "1st Move"
$ws.Cells.Item(60, 2).Value2 = 21617
$ws.Cells.Item(60, 3).Value2 = 22007
$ws.Cells.Item(60, 5).Value2 = 21615

# --- D column: extend the IF(B>0,C-B,0) formula down to row 60 ---
# (D39:D59 is already an existing shared formula group; only add D60)
$ws.Cells.Item(60, 4).Formula = "=IF(B60 >  0,C60-B60, 0)"

# --- F column: new shared formula group covering F58:F60 ---
$ws.Range("F58:F60").Formula = "=IF(B58 >  0,E58-B58, 0)"
